# Auto-generated Excel COM-interop script
# Applies numeric/string cell updates per the target diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = [double]"0"
$ws.Range("F2").Value = [double]"24.81000000000044"
$ws.Range("H2").Value = [double]"0.0774111501817405"
$ws.Range("I2").Value = [double]"0.0774111501817405"
$ws.Range("L2").Value = [double]"26.18657039494436"
$ws.Range("M2").Value = "[-1.5028629946701457, 53.87600378455886]"
$ws.Range("N2").Value = [double]"0.06321126537859478"
$ws.Range("O2").Value = [double]"0.06321126537859478"
$ws.Range("P2").Value = [double]"1.641552918091965"
$ws.Range("Q2").Value = "[0.19497371824080822, 3.088132117943122]"
$ws.Range("R2").Value = [double]"0.02704218021045057"
$ws.Range("S2").Value = [double]"0.02704218021045057"
$ws.Range("T2").Value = [double]"58.5368055023941"
$ws.Range("U2").Value = "[42.24209214617926, 74.83151885860894]"
$ws.Range("V2").Value = [double]"4.566565348085305e-09"
$ws.Range("W2").Value = [double]"4.566565348085305e-09"
$ws.Range("X2").Value = [double]"18.32810810810843"
$ws.Range("Y2").Value = [double]"12.61609609609631"
$ws.Range("Z2").Value = [double]"24.04012012012054"
$ws.Range("F3").Value = [double]"24.81000000000044"
$ws.Range("H3").Value = [double]"7.496604663759232e-05"
$ws.Range("I3").Value = [double]"7.496604663759232e-05"
$ws.Range("L3").Value = [double]"50.33688176759924"
$ws.Range("M3").Value = "[26.50207396996946, 74.17168956522902]"
$ws.Range("N3").Value = [double]"0.0001050101078436949"
$ws.Range("O3").Value = [double]"0.0001050101078436949"
$ws.Range("P3").Value = [double]"1.515763422452734"
$ws.Range("Q3").Value = "[0.9119738433844269, 2.119553001521041]"
$ws.Range("R3").Value = [double]"7.626395422066778e-06"
$ws.Range("S3").Value = [double]"7.626395422066778e-06"
$ws.Range("T3").Value = [double]"69.12284166615062"
$ws.Range("U3").Value = "[54.14243140957558, 84.10325192272566]"
$ws.Range("V3").Value = [double]"4.902966921349616e-12"
$ws.Range("W3").Value = [double]"4.902966921349616e-12"
$ws.Range("X3").Value = [double]"18.82480480480513"
$ws.Range("Y3").Value = [double]"16.44066066066095"
$ws.Range("Z3").Value = [double]"21.20894894894931"
$ws.Range("F4").Value = [double]"24.81000000000044"
$ws.Range("H4").Value = [double]"1.728867370376364e-07"
$ws.Range("I4").Value = [double]"1.728867370376364e-07"
$ws.Range("L4").Value = [double]"58.07023966377325"
$ws.Range("M4").Value = "[35.07268138619692, 81.06779794134958]"
$ws.Range("N4").Value = [double]"6.910934175152761e-06"
$ws.Range("O4").Value = [double]"6.910934175152761e-06"
$ws.Range("P4").Value = [double]"1.150973885098963"
$ws.Range("Q4").Value = "[0.7484474990534249, 1.553500271144502]"
$ws.Range("R4").Value = [double]"7.105039623311882e-07"
$ws.Range("S4").Value = [double]"7.105039623311882e-07"
$ws.Range("T4").Value = [double]"67.20653035517748"
$ws.Range("U4").Value = "[54.79521770365361, 79.61784300670135]"
$ws.Range("V4").Value = [double]"3.197442310920451e-14"
$ws.Range("W4").Value = [double]"3.197442310920451e-14"
$ws.Range("X4").Value = [double]"20.26522522522558"
$ws.Range("Y4").Value = [double]"18.67579579579613"
$ws.Range("Z4").Value = [double]"21.85465465465504"
$ws.Range("B5").Value = [double]"0"
$ws.Range("F5").Value = [double]"24.81000000000044"
$ws.Range("H5").Value = [double]"0.07373815744808077"
$ws.Range("I5").Value = [double]"0.07373815744808077"
$ws.Range("L5").Value = [double]"28.58766830932422"
$ws.Range("M5").Value = "[-1.9494910441424835, 59.12482766279092]"
$ws.Range("N5").Value = [double]"0.0658249627021581"
$ws.Range("O5").Value = [double]"0.0658249627021581"
$ws.Range("P5").Value = [double]"1.779921363295118"
$ws.Range("Q5").Value = "[0.22013161736865516, 3.3397111092215814]"
$ws.Range("R5").Value = [double]"0.02624180334097548"
$ws.Range("S5").Value = [double]"0.02624180334097548"
$ws.Range("T5").Value = [double]"70.97488018016639"
$ws.Range("U5").Value = "[53.35458056907128, 88.5951797912615]"
$ws.Range("V5").Value = [double]"2.363034212748971e-10"
$ws.Range("W5").Value = [double]"2.363034212748971e-10"
$ws.Range("X5").Value = [double]"17.78174174174206"
$ws.Range("Y5").Value = [double]"11.62270270270291"
$ws.Range("Z5").Value = [double]"23.9407807807812"
$ws.Range("B6").Value = [double]"1"
$ws.Range("F6").Value = [double]"24.81000000000044"
$ws.Range("H6").Value = [double]"6.93062969473468e-05"
$ws.Range("I6").Value = [double]"6.93062969473468e-05"
$ws.Range("L6").Value = [double]"39.59059971082596"
$ws.Range("M6").Value = "[17.271334925481145, 61.90986449617078]"
$ws.Range("N6").Value = [double]"0.00085626945144468"
$ws.Range("O6").Value = [double]"0.00085626945144468"
$ws.Range("P6").Value = [double]"0.823921196436963"
$ws.Range("Q6").Value = "[0.2956053147521933, 1.3522370781217328]"
$ws.Range("R6").Value = [double]"0.002973583628305798"
$ws.Range("S6").Value = [double]"0.002973583628305798"
$ws.Range("T6").Value = [double]"58.76931884316395"
$ws.Range("U6").Value = "[47.27989393012521, 70.25874375620269]"
$ws.Range("V6").Value = [double]"2.029487689014786e-13"
$ws.Range("W6").Value = [double]"2.029487689014786e-13"
$ws.Range("X6").Value = [double]"21.55663663663702"
$ws.Range("Y6").Value = [double]"19.47051051051086"
$ws.Range("Z6").Value = [double]"23.64276276276318"
$ws.Range("B7").Value = [double]"1"
$ws.Range("F7").Value = [double]"24.81000000000044"
$ws.Range("H7").Value = [double]"6.836269303500675e-05"
$ws.Range("I7").Value = [double]"6.836269303500675e-05"
$ws.Range("L7").Value = [double]"42.49611108491517"
$ws.Range("M7").Value = "[18.43325865230166, 66.55896351752868]"
$ws.Range("N7").Value = [double]"0.0008970376546739711"
$ws.Range("O7").Value = [double]"0.0008970376546739711"
$ws.Range("P7").Value = [double]"0.5849211547224238"
$ws.Range("Q7").Value = "[0.03144737390980712, 1.1383949355350405]"
$ws.Range("R7").Value = [double]"0.03880155801719476"
$ws.Range("S7").Value = [double]"0.03880155801719476"
$ws.Range("T7").Value = [double]"61.68249535927101"
$ws.Range("U7").Value = "[49.44966659525615, 73.91532412328587]"
$ws.Range("V7").Value = [double]"3.197442310920451e-13"
$ws.Range("W7").Value = [double]"3.197442310920451e-13"
$ws.Range("X7").Value = [double]"22.50036036036076"
$ws.Range("Y7").Value = [double]"20.31489489489525"
$ws.Range("Z7").Value = [double]"24.68582582582626"
$ws.Range("F8").Value = [double]"22.90000000000014"
$ws.Range("H8").Value = [double]"0.001904207304858963"
$ws.Range("I8").Value = [double]"0.001904207304858963"
$ws.Range("L8").Value = [double]"33.1790150422916"
$ws.Range("M8").Value = "[11.347063501067204, 55.01096658351599]"
$ws.Range("N8").Value = [double]"0.003714030427353965"
$ws.Range("O8").Value = [double]"0.003714030427353965"
$ws.Range("P8").Value = [double]"-0.2390000417145384"
$ws.Range("Q8").Value = "[-1.0629212381515005, 0.5849211547224238]"
$ws.Range("R8").Value = [double]"0.5619739861429405"
$ws.Range("S8").Value = [double]"0.5619739861429405"
$ws.Range("T8").Value = [double]"42.69011269362457"
$ws.Range("U8").Value = "[30.005543576949904, 55.374681810299236]"
$ws.Range("V8").Value = [double]"2.172604474104389e-08"
$ws.Range("W8").Value = [double]"2.172604474104389e-08"
$ws.Range("X8").Value = [double]"0.8710710710710714"
$ws.Range("Y8").Value = [double]"-2.131831831831848"
$ws.Range("Z8").Value = [double]"3.873973973973991"
$ws.Range("B9").Value = [double]"1"
$ws.Range("F9").Value = [double]"22.90000000000014"
$ws.Range("H9").Value = [double]"0.003161210226797384"
$ws.Range("I9").Value = [double]"0.003161210226797384"
$ws.Range("L9").Value = [double]"37.8704125746484"
$ws.Range("M9").Value = "[10.09872723224747, 65.64209791704933]"
$ws.Range("N9").Value = [double]"0.008630217293346787"
$ws.Range("O9").Value = [double]"0.008630217293346787"
$ws.Range("P9").Value = [double]"0.9748685912040402"
$ws.Range("Q9").Value = "[0.24528951649650121, 1.7044476659115793]"
$ws.Range("R9").Value = [double]"0.009957010587807069"
$ws.Range("S9").Value = [double]"0.009957010587807069"
$ws.Range("T9").Value = [double]"55.04351983866211"
$ws.Range("U9").Value = "[40.32046773568669, 69.76657194163752]"
$ws.Range("V9").Value = [double]"1.68006897283135e-09"
$ws.Range("W9").Value = [double]"1.68006897283135e-09"
$ws.Range("X9").Value = [double]"19.34694694694706"
$ws.Range("Y9").Value = [double]"16.68788788788799"
$ws.Range("Z9").Value = [double]"22.00600600600614"
$ws.Range("F10").Value = [double]"22.90000000000014"
$ws.Range("H10").Value = [double]"0.0007101407974131613"
$ws.Range("I10").Value = [double]"0.0007101407974131613"
$ws.Range("L10").Value = [double]"50.31054654170806"
$ws.Range("M10").Value = "[16.239458173173418, 84.38163491024271]"
$ws.Range("N10").Value = [double]"0.004710306284302135"
$ws.Range("O10").Value = [double]"0.004710306284302135"
$ws.Range("P10").Value = [double]"0.7736053981812709"
$ws.Range("Q10").Value = "[0.1823947686768852, 1.3648160276856567]"
$ws.Range("R10").Value = [double]"0.01148516485107143"
$ws.Range("S10").Value = [double]"0.01148516485107143"
$ws.Range("T10").Value = [double]"68.9706482839203"
$ws.Range("U10").Value = "[51.697482558429215, 86.24381400941138]"
$ws.Range("V10").Value = [double]"2.992888159525364e-10"
$ws.Range("W10").Value = [double]"2.992888159525364e-10"
$ws.Range("X10").Value = [double]"20.0804804804806"
$ws.Range("Y10").Value = [double]"17.92572572572583"
$ws.Range("Z10").Value = [double]"22.23523523523537"
